# Applies the GDD.docx edits described in the commit:
# "Reworked AAL and corrected errors in AAL and GDD"
#
# 1. "Magic Attack" (Attacks heading)      -> italicize "Magic" only.
# 2. "Fires a large magical ball..." body  -> italicize "magical" only.
# 3. "Magic Missile (" (Spells heading)    -> italicize "Magic" only.
# 4. Move the "_GoBack" bookmark from the end of the "Arc Lightning" bullet
#    to the middle of the "Bubble Trap" bullet (right after "traps ").

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $expectedSubstring, $hintIndex) {
    # Prefer the hinted paragraph index (fast path); fall back to a linear
    # scan of the whole document if the hint doesn't line up, so the script
    # keeps working even if paragraph numbering ever shifts.
    if ($hintIndex -ge 1 -and $hintIndex -le $doc.Paragraphs.Count) {
        $candidate = $doc.Paragraphs.Item($hintIndex)
        if ($candidate.Range.Text -like "*$expectedSubstring*") {
            return $candidate
        }
    }
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        if ($candidate.Range.Text -like "*$expectedSubstring*") {
            return $candidate
        }
    }
    return $null
}

# --- 1. "Magic Attack" heading: italicize "Magic" -------------------------
$p = Get-ParagraphByText $d "Magic Attack" 43
if ($p -ne $null) {
    $rng = $p.Range
    $found = $rng.Find.Execute("Magic", $false, $true, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Italic = 1
    }
}

# --- 2. "Fires a large magical ball..." body: italicize "magical" --------
$p = Get-ParagraphByText $d "Fires a large magical ball" 44
if ($p -ne $null) {
    $rng = $p.Range
    $found = $rng.Find.Execute("magical", $false, $true, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Italic = 1
    }
}

# --- 3. "Magic Missile (" heading: italicize "Magic" ----------------------
$p = Get-ParagraphByText $d "Magic Missile (" 108
if ($p -ne $null) {
    $rng = $p.Range
    $found = $rng.Find.Execute("Magic", $false, $true, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Italic = 1
    }
}

# --- 4. Move the "_GoBack" bookmark into the "Bubble Trap" bullet ---------
# Adding a bookmark named "_GoBack" redefines/relocates the existing one
# (bookmark names are unique), so this both removes it from the end of the
# "Arc Lightning" bullet and inserts it between "traps " and "enemies..."
# in the "Bubble Trap" bullet, splitting that run in two in the process.
$p = Get-ParagraphByText $d "Bubble Trap (Low fire rate" 114
if ($p -ne $null) {
    $rng = $p.Range
    $found = $rng.Find.Execute("traps ", $false, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $splitPoint = $rng.Duplicate
        $splitPoint.Collapse(0)
        $d.Bookmarks.Add("_GoBack", $splitPoint)
    }
}
